$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualizacao de bases das ligas, do dia: 27-03-2024 as 20:23
# Rows 9 <-> 10 : match data swapped between the two adjacent rows (id column A untouched)
$ws.Range('B9').Value = 6865281
$ws.Range('F9').Value = 'GOSK Gabela'
$ws.Range('G9').Value = 'Zvijezda 09'
$ws.Range('H9').Value = 2
$ws.Range('K9').Value = 1.75
$ws.Range('L9').Value = 4
$ws.Range('M9').Value = 3.5
$ws.Range('N9').Value = 1.75
$ws.Range('O9').Value = 4
$ws.Range('R9').Value = 1.8
$ws.Range('S9').Value = 2
$ws.Range('W9').Value = 0.75
$ws.Range('Z9').Value = 0.8

$ws.Range('B10').Value = 6865285
$ws.Range('F10').Value = 'NK Igman Konjic'
$ws.Range('G10').Value = 'Sloga'
$ws.Range('H10').Value = 1
$ws.Range('K10').Value = 2
$ws.Range('L10').Value = 3.4
$ws.Range('M10').Value = 3.2
$ws.Range('N10').Value = 1.909
$ws.Range('O10').Value = 3.5
$ws.Range('R10').Value = 1.95
$ws.Range('S10').Value = 1.85
$ws.Range('W10').Value = 0.909
$ws.Range('Z10').Value = 0.95

# Rows 29 <-> 30 : match data swapped between the two adjacent rows
$ws.Range('B29').Value = 6865296
$ws.Range('F29').Value = 'Velez Mostar'
$ws.Range('G29').Value = 'Zeljeznicar'
$ws.Range('H29').Value = 1
$ws.Range('I29').Value = 0
$ws.Range('K29').Value = 1.909
$ws.Range('L29').Value = 3.2
$ws.Range('M29').Value = 3.6
$ws.Range('N29').Value = 1.95
$ws.Range('O29').Value = 3.2
$ws.Range('P29').Value = 3.4
$ws.Range('Q29').Value = -0.5
$ws.Range('R29').Value = 2.025
$ws.Range('S29').Value = 1.775
$ws.Range('T29').Value = 2.25
$ws.Range('U29').Value = 1.9
$ws.Range('V29').Value = 1.9
$ws.Range('W29').Value = 0.95
$ws.Range('Z29').Value = 1.025
$ws.Range('AB29').Value = -1
$ws.Range('AC29').Value = 0.8999999999999999

$ws.Range('B30').Value = 6865295
$ws.Range('F30').Value = 'FK Tuzla City'
$ws.Range('G30').Value = 'NK Igman Konjic'
$ws.Range('H30').Value = 3
$ws.Range('I30').Value = 1
$ws.Range('K30').Value = 1.8
$ws.Range('L30').Value = 3.4
$ws.Range('M30').Value = 3.8
$ws.Range('N30').Value = 1.615
$ws.Range('O30').Value = 3.5
$ws.Range('P30').Value = 4.5
$ws.Range('Q30').Value = -0.75
$ws.Range('R30').Value = 1.85
$ws.Range('S30').Value = 1.95
$ws.Range('T30').Value = 2.75
$ws.Range('U30').Value = 2
$ws.Range('V30').Value = 1.8
$ws.Range('W30').Value = 0.615
$ws.Range('Z30').Value = 0.8500000000000001
$ws.Range('AB30').Value = 1
$ws.Range('AC30').Value = -1

# Rows 76 <-> 77 : match data swapped between the two adjacent rows
$ws.Range('B76').Value = 6865328
$ws.Range('F76').Value = 'Siroki Brijeg'
$ws.Range('G76').Value = 'NK Posusje'
$ws.Range('H76').Value = 1
$ws.Range('J76').Value = 'D'
$ws.Range('K76').Value = 2
$ws.Range('L76').Value = 3
$ws.Range('M76').Value = 3.5
$ws.Range('N76').Value = 2.1
$ws.Range('O76').Value = 3
$ws.Range('P76').Value = 3.3
$ws.Range('Q76').Value = -0.25
$ws.Range('R76').Value = 1.825
$ws.Range('S76').Value = 1.975
$ws.Range('T76').Value = 2
$ws.Range('U76').Value = 1.825
$ws.Range('V76').Value = 1.975
$ws.Range('W76').Value = -1
$ws.Range('X76').Value = 2
$ws.Range('Z76').Value = -0.5
$ws.Range('AA76').Value = 0.4875
$ws.Range('AB76').Value = 0
$ws.Range('AC76').Value = -0

$ws.Range('B77').Value = 6865377
$ws.Range('F77').Value = 'Zrinjski Mostar'
$ws.Range('G77').Value = 'FK Tuzla City'
$ws.Range('H77').Value = 3
$ws.Range('J77').Value = 'H'
$ws.Range('K77').Value = 1.333
$ws.Range('L77').Value = 5
$ws.Range('M77').Value = 6
$ws.Range('N77').Value = 1.166
$ws.Range('O77').Value = 6.5
$ws.Range('P77').Value = 13
$ws.Range('Q77').Value = -2
$ws.Range('R77').Value = 1.9
$ws.Range('S77').Value = 1.9
$ws.Range('T77').Value = 3.25
$ws.Range('U77').Value = 1.95
$ws.Range('V77').Value = 1.85
$ws.Range('W77').Value = 0.1659999999999999
$ws.Range('X77').Value = -1
$ws.Range('Z77').Value = 0
$ws.Range('AA77').Value = -0
$ws.Range('AB77').Value = 0.95
$ws.Range('AC77').Value = -1

# Rows 87 <-> 88 : match data swapped between the two adjacent rows
$ws.Range('B87').Value = 7505497
$ws.Range('F87').Value = 'Zeljeznicar'
$ws.Range('G87').Value = 'NK Posusje'
$ws.Range('I87').Value = 1
$ws.Range('J87').Value = 'D'
$ws.Range('K87').Value = 1.65
$ws.Range('L87').Value = 3.4
$ws.Range('M87').Value = 4.75
$ws.Range('N87').Value = 1.8
$ws.Range('O87').Value = 3.2
$ws.Range('P87').Value = 4.2
$ws.Range('Q87').Value = -0.5
$ws.Range('R87').Value = 1.825
$ws.Range('S87').Value = 1.975
$ws.Range('T87').Value = 2
$ws.Range('U87').Value = 1.75
$ws.Range('V87').Value = 2.05
$ws.Range('W87').Value = -1
$ws.Range('X87').Value = 2.2
$ws.Range('Z87').Value = -1
$ws.Range('AA87').Value = 0.9750000000000001
$ws.Range('AB87').Value = 0
$ws.Range('AC87').Value = -0

$ws.Range('B88').Value = 7505495
$ws.Range('F88').Value = 'Sloga'
$ws.Range('G88').Value = 'Zvijezda 09'
$ws.Range('I88').Value = 0
$ws.Range('J88').Value = 'H'
$ws.Range('K88').Value = 1.444
$ws.Range('L88').Value = 4.2
$ws.Range('M88').Value = 5.5
$ws.Range('N88').Value = 1.5
$ws.Range('O88').Value = 4.2
$ws.Range('P88').Value = 5.25
$ws.Range('Q88').Value = -1
$ws.Range('R88').Value = 1.8
$ws.Range('S88').Value = 2
$ws.Range('T88').Value = 2.75
$ws.Range('U88').Value = 1.775
$ws.Range('V88').Value = 2.025
$ws.Range('W88').Value = 0.5
$ws.Range('X88').Value = -1
$ws.Range('Z88').Value = 0
$ws.Range('AA88').Value = -0
$ws.Range('AB88').Value = -1
$ws.Range('AC88').Value = 1.025

# Rows 99 <-> 100 : match data swapped between the two adjacent rows
$ws.Range('B99').Value = 6865343
$ws.Range('F99').Value = 'Sloga'
$ws.Range('G99').Value = 'NK Posusje'
$ws.Range('I99').Value = 0
$ws.Range('J99').Value = 'H'
$ws.Range('K99').Value = 1.909
$ws.Range('L99').Value = 3.3
$ws.Range('M99').Value = 3.5
$ws.Range('N99').Value = 2.2
$ws.Range('O99').Value = 2.8
$ws.Range('P99').Value = 3.3
$ws.Range('Q99').Value = -0.25
$ws.Range('R99').Value = 1.95
$ws.Range('S99').Value = 1.85
$ws.Range('T99').Value = 1.75
$ws.Range('U99').Value = 1.875
$ws.Range('V99').Value = 1.925
$ws.Range('W99').Value = 1.2
$ws.Range('Y99').Value = -1
$ws.Range('Z99').Value = 0.95
$ws.Range('AB99').Value = -1
$ws.Range('AC99').Value = 0.925

$ws.Range('B100').Value = 6864639
$ws.Range('F100').Value = 'Zvijezda 09'
$ws.Range('G100').Value = 'Borac Banja Luka'
$ws.Range('I100').Value = 2
$ws.Range('J100').Value = 'A'
$ws.Range('K100').Value = 11
$ws.Range('L100').Value = 6
$ws.Range('M100').Value = 1.2
$ws.Range('N100').Value = 10
$ws.Range('O100').Value = 6.5
$ws.Range('P100').Value = 1.181
$ws.Range('Q100').Value = 2
$ws.Range('R100').Value = 1.825
$ws.Range('S100').Value = 1.975
$ws.Range('T100').Value = 3
$ws.Range('U100').Value = 1.9
$ws.Range('V100').Value = 1.9
$ws.Range('W100').Value = -1
$ws.Range('Y100').Value = 0.181
$ws.Range('Z100').Value = 0.825
$ws.Range('AB100').Value = 0
$ws.Range('AC100').Value = -0

# Rows 137-139 : fill in previously-missing FTHG/FTAG/FTR and recompute the
# profit/loss (PL*) columns that depended on them
$ws.Range('H137').Value = 1
$ws.Range('I137').Value = 0
$ws.Range('J137').Value = 'H'
$ws.Range('W137').Value = 0.6499999999999999
$ws.Range('X137').Value = -1
$ws.Range('Y137').Value = -1
$ws.Range('Z137').Value = 0.45
$ws.Range('AA137').Value = -0.5
$ws.Range('AB137').Value = -1
$ws.Range('AC137').Value = 1.025

$ws.Range('H138').Value = 3
$ws.Range('I138').Value = 1
$ws.Range('J138').Value = 'H'
$ws.Range('N138').Value = 1.25
$ws.Range('O138').Value = 4.5
$ws.Range('P138').Value = 10
$ws.Range('Q138').Value = -1.5
$ws.Range('R138').Value = 1.85
$ws.Range('S138').Value = 1.95
$ws.Range('T138').Value = 2.5
$ws.Range('U138').Value = 1.825
$ws.Range('V138').Value = 1.975
$ws.Range('W138').Value = 0.25
$ws.Range('X138').Value = -1
$ws.Range('Y138').Value = -1
$ws.Range('Z138').Value = 0.8500000000000001
$ws.Range('AA138').Value = -1
$ws.Range('AB138').Value = 0.825
$ws.Range('AC138').Value = -1

$ws.Range('H139').Value = 2
$ws.Range('I139').Value = 0
$ws.Range('J139').Value = 'H'
$ws.Range('W139').Value = 0.3
$ws.Range('X139').Value = -1
$ws.Range('Y139').Value = -1
$ws.Range('Z139').Value = 0.825
$ws.Range('AA139').Value = -1
$ws.Range('AB139').Value = -0.5
$ws.Range('AC139').Value = 0.5125
